$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 387, pushing existing rows 387:410 down to 388:411
$ws.Rows.Item(387).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A387").Value = 4
$ws.Range("B387").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C387").Value = "Los Lagos"
$ws.Range("D387").Value = 44746
$ws.Range("E387").Value = 10
$ws.Range("F387").Value = 100114001
$ws.Range("G387").Value = "Papa"
$ws.Range("H387").Value = "Patagonia"
$ws.Range("I387").Value = "1a (guarda)"
$ws.Range("J387").Value = 300
$ws.Range("K387").Value = 7000
$ws.Range("L387").Value = 8000
$ws.Range("M387").Value = 7500
$ws.Range("N387").Value = "$/saco 25 kilos"
$ws.Range("O387").Value = "Provincia de Llanquihue"
$ws.Range("P387").Value = 300
$ws.Range("Q387").Value = 25
$ws.Range("R387").Value = "Hortaliza"
